$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.065.17"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.756.89"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.33"
$ws.Range("E5").Value = "  -0.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.03"
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
$ws.Range("E8").Value = "  -2.64%  "

# Row 9
$ws.Range("E9").Value = "  -1.87%  "

# Row 10
$ws.Range("E10").Value = "  +3.85%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.85"
$ws.Range("E11").Value = "  -14.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.246.99"
$ws.Range("E13").Value = "  +0.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.02"
$ws.Range("E14").Value = "  -2.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.690.33"
$ws.Range("E15").Value = "  +0.05%  "

# Row 16
$ws.Range("E16").Value = "  -2.99%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.762.07"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.22"
$ws.Range("E18").Value = "  +0.09%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.86"
$ws.Range("E19").Value = "  -2.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.29"
$ws.Range("E20").Value = "  -1.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").Value = "  -3.97%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23
$ws.Range("E23").Value = "  -1.08%  "

# Row 24
$ws.Range("E24").Value = "  -1.65%  "

# Row 25
$ws.Range("E25").Value = "  -0.99%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.62"
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").Value = "  -1.34%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  +1.32%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.95"
$ws.Range("E30").Value = "  -3.86%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.25"
$ws.Range("E31").Value = "  -1.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.08"
$ws.Range("E32").Value = "  -2.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.96"
$ws.Range("E33").Value = "  -0.54%  "

# Row 34
$ws.Range("E34").Value = "  -2.17%  "

# Row 35
$ws.Range("E35").Value = "  +1.67%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("E37").Value = "  -0.75%  "

# Row 38
$ws.Range("E38").Value = "  -1.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "350.28"
$ws.Range("E39").Value = "  +3.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.30"
$ws.Range("E40").Value = "  +1.08%  "

# Row 41
$ws.Range("E41").Value = "  -1.39%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.26"
$ws.Range("E42").Value = "  -0.79%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.56"
$ws.Range("E43").Value = "  -1.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.92"
$ws.Range("E44").Value = "  -2.39%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0590"
$ws.Range("E45").Value = "  -2.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "137.07"
$ws.Range("E46").Value = "  -0.27%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.633"
$ws.Range("E47").Value = "  -1.92%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0255"
$ws.Range("E48").Value = "  -1.71%  "

# Row 49
$ws.Range("E49").Value = "  -0.37%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.22%  "

# Row 51
$ws.Range("E51").Value = "  -0.01%  "

